# Applies the changes described by the commit:
# "Added Pop up alerting problem on connection with network or broker MQTT
#  and credential invalid. Also I added function to parse response of
#  credential and storage liters available."
#
# Concretely, on the "Translation" worksheet:
#   - Fixes the spelling of the existing "Invalid credencial, / try again
#     please" text used in row 33 (F33) to "Invalid credential, / try
#     again please".
#   - Adds three new rows (39, 40, 41) with new SingleUseId text entries
#     for: an "OK" button label, a "network connection failed" message,
#     and an "MQTT server problem" message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Fix typo: "credencial" -> "credential" in the existing popup text (row 33)
$ws.Range("F33").Value = "Invalid credential, `ntry again please"
$ws.Rows.Item(33).AutoFit() | Out-Null

# New row 39: OK button for the new pop-up(s)
$ws.Range("B39").Value = "SingleUseId34"
$ws.Range("C39").Value = "Typography_label"
$ws.Range("D39").Value = "Center"
$ws.Range("E39").Value = "LTR"
$ws.Range("F39").Value = "OK"

# New row 40: "Couldn't connect to network" pop-up text
$ws.Range("B40").Value = "SingleUseId35"
$ws.Range("C40").Value = "Typography_label"
$ws.Range("D40").Value = "Center"
$ws.Range("E40").Value = "LTR"
$ws.Range("F40").Value = "Couldn't connect to network,`ntry again please."

# New row 41: "Problem with MQTT server" pop-up text
$ws.Range("B41").Value = "SingleUseId36"
$ws.Range("C41").Value = "Typography_label"
$ws.Range("D41").Value = "Center"
$ws.Range("E41").Value = "LTR"
$ws.Range("F41").Value = "Problem with MQTT`nserver"

# Keep row heights at their natural auto-fit value (avoid leaving an
# explicit custom row height behind just because the new text wraps
# across multiple lines).
$ws.Rows.Item(39).AutoFit() | Out-Null
$ws.Rows.Item(40).AutoFit() | Out-Null
$ws.Rows.Item(41).AutoFit() | Out-Null
